# "Version finale enregistrement excel"
# Extend the logged sample row from 13 columns (B:N) to 32 columns (B:AG),
# relabel the header row with the new running sample index, and rewrite
# the data row with the "Setup done" marker plus the "0.00" formatted
# readings produced by the final python/arduino logger.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Row 1: header / sample-index row -------------------------------
# B1:N1 already hold 0..12 using style index 1 (bold font, boxed border,
# centered). Grow the row out to column AG (new max value 31), copying
# N1's formatting across the new cells first so no extra style entries
# get created, then overwrite each cell's value with the next index.
$ws.Range("N1").Copy($ws.Range("O1:AG1"))

$headerCols = @("O","P","Q","R","S","T","U","V","W","X","Y","Z","AA","AB","AC","AD","AE","AF","AG")
$headerVal = 13
foreach ($col in $headerCols) {
    $ws.Range($col + "1").Value = $headerVal
    $headerVal = $headerVal + 1
}

# ---- Row 2: data row --------------------------------------------------
# B2 becomes the "Setup done" marker written once at logger start-up.
# C2:AF2 become the "0.00" placeholder reading written on every tick.
# AG2 is a trailing empty text cell (mirrors the previous trailing empty
# cell that used to sit at N2).
#
# A leading apostrophe forces Excel to store digit-looking / empty text
# ("0.00", "") as literal text instead of auto-parsing it as a number;
# resetting the range back to the "Normal" style afterwards drops the
# number-format bookkeeping that the text-entry added, so the cells keep
# the workbook's default (unstyled) look, same as before the edit.

$ws.Range("B2").Value = "'Setup done`n"

$dataCols = @("C","D","E","F","G","H","I","J","K","L","M","N","O","P","Q","R","S","T","U","V","W","X","Y","Z","AA","AB","AC","AD","AE","AF")
foreach ($col in $dataCols) {
    $ws.Range($col + "2").Value = "'0.00`n"
}

$ws.Range("AG2").Value = "'"

$ws.Range("B2:AG2").Style = "Normal"
